{"js": "// Apply four small textual corrections to paragraphs in the document body.\n// Each target string is the sole text of a single run in its own paragraph,\n// so a search + Replace insertText keeps the existing run formatting intact.\nconst replacements = [\n  { find: \"\u042d\u0442\u043e \u0442\u043e \u0447\u0442\u043e \u043d\u0430 \u0433\u0440\u0430\u043d\u0438 \u043f\u0440\u0438\u0440\u043e\u0434\u044b\", replace: \"\u042d\u0442\u043e \u0442\u043e, \u0447\u0442\u043e \u043d\u0430 \u0433\u0440\u0430\u043d\u0438 \u043f\u0440\u0438\u0440\u043e\u0434\u044b\" },\n  { find: \"\u041b\u0435\u0431\u0435\u0437\u0438\u0442\", replace: \"\u0421\u043a\u0430\u0431\u0440\u0435\u0437\u0438\u0442\" },\n  { find: \"\u041f\u0440\u0430\u0434\u0435\u0434\u044b \u043d\u0430\u0448\u0438 \u0431\u044b\u043b\u0438\", replace: \"\u0418 \u043f\u0440\u0430\u0434\u0435\u0434\u044b - \u0431\u044b\u043b\u0438\" },\n  { find: \"\u0425\u0440\u0430\u043d\u044e \u0447\u0435\u0441\u0442\u044c-\u043d\u0430\u0434\u0435\u0436\u0434\u0443\", replace: \"\u0425\u0440\u0430\u043d\u044e \u0447\u0435\u0441\u0442\u044c - \u043d\u0430\u0434\u0435\u0436\u0434\u0443\" },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  results.items[0].insertText(replace, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Apply four small textual corrections to paragraphs in the document body.\n# Each target string is the sole text of a single run in its own paragraph,\n# so a Find/Replace over the whole document body reliably hits exactly one\n# location for every pair below.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"\u042d\u0442\u043e \u0442\u043e \u0447\u0442\u043e \u043d\u0430 \u0433\u0440\u0430\u043d\u0438 \u043f\u0440\u0438\u0440\u043e\u0434\u044b\"; Replace = \"\u042d\u0442\u043e \u0442\u043e, \u0447\u0442\u043e \u043d\u0430 \u0433\u0440\u0430\u043d\u0438 \u043f\u0440\u0438\u0440\u043e\u0434\u044b\" },\n    @{ Find = \"\u041b\u0435\u0431\u0435\u0437\u0438\u0442\"; Replace = \"\u0421\u043a\u0430\u0431\u0440\u0435\u0437\u0438\u0442\" },\n    @{ Find = \"\u041f\u0440\u0430\u0434\u0435\u0434\u044b \u043d\u0430\u0448\u0438 \u0431\u044b\u043b\u0438\"; Replace = \"\u0418 \u043f\u0440\u0430\u0434\u0435\u0434\u044b - \u0431\u044b\u043b\u0438\" },\n    @{ Find = \"\u0425\u0440\u0430\u043d\u044e \u0447\u0435\u0441\u0442\u044c-\u043d\u0430\u0434\u0435\u0436\u0434\u0443\"; Replace = \"\u0425\u0440\u0430\u043d\u044e \u0447\u0435\u0441\u0442\u044c - \u043d\u0430\u0434\u0435\u0436\u0434\u0443\" }\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair.Find\n    $replaceText = $pair.Replace\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $ok = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n\n    if (-not $ok) {\n        throw \"Text not found: $findText\"\n    }\n}\n"}
